$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D5: status set to "hold" (new distinct text entry + its own style)
$d5 = $ws.Range("D5")
$d5.Value = "hold"
$d5.Font.Color = 1
$d5.Characters(1, 3).Font.Name = "Calibri"

# Helper to mark a bug row as resolved ("已解决") with a resolution date,
# copying the date/number format already used elsewhere in column F.
function Set-Resolved($rowNum, $serial) {
    $ws.Range("D$rowNum").Value = "已解决"

    $ws.Range("F6").Copy() | Out-Null
    $ws.Range("F$rowNum").PasteSpecial(-4122) | Out-Null
    $ws.Range("F$rowNum").Value = $serial
}

Set-Resolved 21 42135
Set-Resolved 22 42135
Set-Resolved 27 42131
Set-Resolved 28 42131
Set-Resolved 29 42135
Set-Resolved 30 42135

$excel.CutCopyMode = $false

# Reflect the scrolled/selected cell from the author's session.
$ws.Range("C30").Select() | Out-Null

Write-Host "Bug list status updated."
